$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...put the player character (P) on a spot..." -> "...($) on..."
#   The run containing "(P) " gets split right after the new "$" so that the
#   saved OOXML has two runs: "...player character ($" and ") ".
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(8)
$f1 = $p1.Range.Duplicate
$f1.Find.Execute("(P)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start1 = $f1.Start   # offset of the "(" character

# Replace the "P" character with "$"
$pChar1 = $d.Range($start1 + 1, $start1 + 2)
$pChar1.Text = "$"

# Force a run break right after the "$" (between "$" and ")") by toggling a
# character formatting property on/off on the remainder of the old run; this
# must be the last edit touching the paragraph so the split sticks.
$tail1 = $d.Range($start1 + 2, $start1 + 4)   # ") " including trailing space
$tail1.Bold = $true
$tail1.Bold = $false

# ---------------------------------------------------------------------------
# Change 2: "You must have a "P" in a stage. ..." -> "...a "$" in a stage..."
#   The run containing " "P" in a stage..." gets split right after the new
#   "$" into: " "$"  and  "" in a stage. This represents ...character."
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(9)
$f2 = $p2.Range.Duplicate
$f2.Find.Execute("P", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $f2.Start   # offset of the "P" character

# Replace the "P" character with "$"
$pChar2 = $d.Range($start2, $start2 + 1)
$pChar2.Text = "$"

$paraEnd2 = $p2.Range.End
$runStart2 = $start2 - 2   # the leading space before the opening smart quote

# First isolate the whole former run (leading space through end of paragraph,
# excluding the paragraph mark) away from the preceding "a" run...
$wholeTail2 = $d.Range($runStart2, $paraEnd2 - 1)
$wholeTail2.Bold = $true
$wholeTail2.Bold = $false

# ...then split that chunk right after the "$" so the closing smart quote and
# the remaining sentence become their own run.
$tail2 = $d.Range($start2 + 1, $paraEnd2 - 1)
$tail2.Bold = $true
$tail2.Bold = $false
